# Fix latency units in report sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header O2: "Utility" -> "Utility (Percent)"
$ws.Range("O2").Value = "Utility (Percent)"

# Append " usec" to the Read Latency columns (I, J, K) for data rows 3-14
foreach ($row in 3..14) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Text
        $cell.Value = "$current usec"
    }
}
